$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2694.0908
$ws.Range("J17").Value = 2694.0908
$ws.Range("L17").Value = 8082.2724
$ws.Range("N17").Value = -8418.2724

$ws.Range("H32").Value = 17721.111
$ws.Range("J32").Value = 17122.5
$ws.Range("L32").Value = 17122.5
$ws.Range("N32").Value = -17774.5

$ws.Range("H74").Value = 4608
$ws.Range("I74").Value = 4608
$ws.Range("K74").Value = 4608
$ws.Range("M74").Value = -3672

$ws.Range("H77").Value = 4608
$ws.Range("I77").Value = 4608
$ws.Range("K77").Value = 23040
$ws.Range("M77").Value = -18360

$ws.Range("H113").Value = 3079.0625
$ws.Range("I113").Value = 3033.9285
$ws.Range("K113").Value = 3033.9285
$ws.Range("M113").Value = 220.0715

$ws.Range("H116").Value = 2746.5
$ws.Range("I116").Value = 2619.75
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2619.75
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 822.25
$ws.Range("N116").Value = -9884

$ws.Range("H121").Value = 1549.7142
$ws.Range("J121").Value = 1549.7142
$ws.Range("L121").Value = 4649.142599999999
$ws.Range("N121").Value = -8143.142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3290.1155
$ws.Range("I32").Value = 2898.24
$ws.Range("K32").Value = 2898.24
$ws.Range("M32").Value = -2611.24

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 20998.4
$ws.Range("J6").Value = 20998.4
$ws.Range("L6").Value = 20998.4
$ws.Range("N6").Value = -21224.4

$ws.Range("H20").Value = 5563
$ws.Range("I20").Value = 1172
$ws.Range("J20").Value = 9954
$ws.Range("K20").Value = 1172
$ws.Range("L20").Value = 9954
$ws.Range("M20").Value = -925
$ws.Range("N20").Value = -10448

$ws.Range("H86").Value = 11620.714
$ws.Range("I86").Value = 20047.834
$ws.Range("J86").Value = 5300.375
$ws.Range("K86").Value = 20047.834
$ws.Range("L86").Value = 5300.375
$ws.Range("M86").Value = -18924.834
$ws.Range("N86").Value = -7546.375

$ws.Range("H89").Value = 11620.714
$ws.Range("I89").Value = 20047.834
$ws.Range("J89").Value = 5300.375
$ws.Range("K89").Value = 100239.17
$ws.Range("L89").Value = 26501.875
$ws.Range("M89").Value = -94623.17
$ws.Range("N89").Value = -37733.875

$ws.Range("H97").Value = 7761.6665
$ws.Range("I97").Value = 7761.6665
$ws.Range("K97").Value = 7761.6665
$ws.Range("M97").Value = -6770.6665

$ws.Range("H99").Value = 1451.4286
$ws.Range("I99").Value = 1043.3334
$ws.Range("J99").Value = 3900
$ws.Range("K99").Value = 1043.3334
$ws.Range("L99").Value = 3900
$ws.Range("M99").Value = 454.6666
$ws.Range("N99").Value = -6896

$ws.Range("H105").Value = 3849.6667
$ws.Range("I105").Value = 3533
$ws.Range("K105").Value = 3533
$ws.Range("M105").Value = -1786

$ws.Range("H134").Value = 1899.6666
$ws.Range("I134").Value = 1899.6666
$ws.Range("K134").Value = 5698.9998
$ws.Range("M134").Value = -3163.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 141.07143
$ws.Range("I7").Value = 85.416664
$ws.Range("J7").Value = 475
$ws.Range("K7").Value = 85.416664
$ws.Range("L7").Value = 475
$ws.Range("M7").Value = 27.583336
$ws.Range("N7").Value = -701

$ws.Range("H16").Value = 496.3
$ws.Range("I16").Value = 536.6667
$ws.Range("J16").Value = 133
$ws.Range("K16").Value = 536.6667
$ws.Range("L16").Value = 133
$ws.Range("M16").Value = -249.6667
$ws.Range("N16").Value = -707

$ws.Range("H22").Value = 3146.6667
$ws.Range("I22").Value = 970
$ws.Range("J22").Value = 7500
$ws.Range("K22").Value = 970
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = -620
$ws.Range("N22").Value = -8200

$ws.Range("H31").Value = 10500.875
$ws.Range("I31").Value = 11573.286
$ws.Range("K31").Value = 11573.286
$ws.Range("M31").Value = -11278.286

$ws.Range("H34").Value = 10500.875
$ws.Range("I34").Value = 11573.286
$ws.Range("K34").Value = 11573.286
$ws.Range("M34").Value = -11371.286

$ws.Range("H58").Value = 3192.8
$ws.Range("I58").Value = 2032.5555
$ws.Range("K58").Value = 2032.5555
$ws.Range("M58").Value = -1829.5555

$ws.Range("H86").Value = 9839.846
$ws.Range("I86").Value = 5265.364
$ws.Range("K86").Value = 5265.364
$ws.Range("M86").Value = -4142.364

$ws.Range("H89").Value = 9839.846
$ws.Range("I89").Value = 5265.364
$ws.Range("K89").Value = 26326.82
$ws.Range("M89").Value = -20710.82

$ws.Range("H105").Value = 1398.8572
$ws.Range("I105").Value = 933.4
$ws.Range("K105").Value = 933.4
$ws.Range("M105").Value = 813.6

$ws.Range("H113").Value = 496.3
$ws.Range("I113").Value = 536.6667
$ws.Range("J113").Value = 133
$ws.Range("K113").Value = 536.6667
$ws.Range("L113").Value = 133
$ws.Range("M113").Value = 1633.3333
$ws.Range("N113").Value = -4473

$ws.Range("H134").Value = 1823.4546
$ws.Range("I134").Value = 1823.4546
$ws.Range("K134").Value = 5470.3638
$ws.Range("M134").Value = -2935.3638

$ws.Range("H136").Value = 3192.8
$ws.Range("I136").Value = 2032.5555
$ws.Range("K136").Value = 6097.666499999999
$ws.Range("M136").Value = -3547.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H34").Value = 8830.833000000001
$ws.Range("J34").Value = 8830.833000000001
$ws.Range("L34").Value = 26492.499
$ws.Range("N34").Value = -26660.499

$ws.Range("H44").Value = 333667.34
$ws.Range("I44").Value = 333667.34
$ws.Range("K44").Value = 1001002.02
$ws.Range("M44").Value = -1000604.02

$ws.Range("H46").Value = 7571.4287
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 10250
$ws.Range("K46").Value = 12000
$ws.Range("L46").Value = 30750
$ws.Range("M46").Value = -11909
$ws.Range("N46").Value = -30932

$ws.Range("H47").Value = 455.66666
$ws.Range("I47").Value = 455.66666
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 1366.99998
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -935.9999800000001
$ws.Range("N47").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50576

$ws.Range("H41").Value = 14750
$ws.Range("I41").Value = 14750
$ws.Range("K41").Value = 14750
$ws.Range("M41").Value = -14395

$ws.Range("H70").Value = 8948
$ws.Range("I70").Value = 7097.6665
$ws.Range("K70").Value = 7097.6665
$ws.Range("M70").Value = -6827.6665

$ws.Range("H73").Value = 8948
$ws.Range("I73").Value = 7097.6665
$ws.Range("K73").Value = 7097.6665
$ws.Range("M73").Value = -6161.6665

$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 25.833334
$ws.Range("I2").Value = 23.636364
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 23.636364
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 88.363636
$ws.Range("N2").Value = -274

$ws.Range("H12").Value = 2668
$ws.Range("J12").Value = 2668
$ws.Range("L12").Value = 2668
$ws.Range("N12").Value = -3008

$ws.Range("H134").Value = 91248.06
$ws.Range("J134").Value = 91248.06
$ws.Range("L134").Value = 91248.06
$ws.Range("N134").Value = -101388.06

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H8").Value = 15000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 15000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 15000
$ws.Range("N8").Value = -15280
$ws.Range("M8").ClearContents()

$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338

$ws.Range("H11").Value = 5003.778
$ws.Range("I11").Value = 5003.778
$ws.Range("K11").Value = 5003.778
$ws.Range("M11").Value = -4861.778

$ws.Range("H13").Value = 18000
$ws.Range("J13").Value = 18000
$ws.Range("L13").Value = 18000
$ws.Range("N13").Value = -18280

$ws.Range("H133").Value = 122500
$ws.Range("J133").Value = 122500
$ws.Range("L133").Value = 122500
$ws.Range("N133").Value = -132620

$ws.Range("H136").Value = 11851.821
$ws.Range("I136").Value = 10718.48
$ws.Range("K136").Value = 32155.44
$ws.Range("M136").Value = -29605.44
